$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Objetivos text (B10/C10) ---
$ws.Range("B10").Value = 'Fomentar a cultura do empreendedorismo; Desenvolver habilidades empreendedoras; Apresentar conhecimentos necessários para a criação de startups. A disciplina é aplicada através de Aprendizagem baseada em Projetos, onde o projeto a ser desenvolvido é da criação de uma startup ao longo do semestre.'
$ws.Range("C10").Value = 'Fomentar a cultura do empreendedorismo; Desenvolver habilidades empreendedoras; Apresentar conhecimentos necessários para a criação de startups. A disciplina é aplicada através de Aprendizagem baseada em Projetos, onde o projeto a ser desenvolvido é da criação de uma startup ao longo do semestre.'

# --- Row 13: drop A13 label; becomes "Docentes responsaveis" value row ---
$ws.Range("A13").Clear()
$ws.Range("B13").Value = '5840560 - Marco Antonio Carvalho Pereira'
$ws.Range("C13").Value = '5840560 - Marco Antonio Carvalho Pereira'
$ws.Rows.Item(13).AutoFit()

# --- Row 14: Programa resumido: label + new summary text ---
$ws.Range("A14").Value = 'Programa resumido:'
$ws.Range("B14").Value = 'Características do Comportamento Empreendedor; Modelo de Negócios; Produto mínimo viável; Plano de Negócios.'
$ws.Range("C14").Value = 'Características do Comportamento Empreendedor; Modelo de Negócios; Produto mínimo viável; Plano de Negócios.'

# --- Row 15: Short syllabus: label + existing EN short text (row height 120 -> 60) ---
$ws.Range("A15").Value = 'Short syllabus:'
$ws.Range("B15").Value = 'Characteristics of Entrepreneurial Behavior. Business Model. Minimum Viable Product. Business Plan.'
$ws.Range("C15").Value = 'Characteristics of Entrepreneurial Behavior. Business Model. Minimum Viable Product. Business Plan.'
$ws.Rows.Item(15).RowHeight = 60

# --- Row 16: Programa: label + new PT full syllabus text ---
$ws.Range("A16").Value = 'Programa:'
$ws.Range("B16").Value = '1.Características do Comportamento Empreendedor: Busca de oportunidades e iniciativa. Correr riscos calculados. Exigência de qualidade e eficiência. Persistência. Comprometimento. Busca de informações. Estabelecimento de metas. Monitoramento e planejamento sistemático. Persuasão e rede contatos. Independência e autoconfiança.2.Modelo de Negócios (Lean Canvas): Problema. Segmento de Clientes. Proposta de Valor Única. Solução. Métricas-Chave. Canais. Estrutura de Custos. Fluxos de Receita. Vantagem Injusta.3.Produto mínimo viável: Ciclo Construir-Mensurar-Aprender. Valor da vida útil do cliente.4.Plano de Negócios: Marketing, Finanças, Recursos Humanos, Desenvolvimento de Produtos e Tecnologia da Informação e Comunicação.'
$ws.Range("C16").Value = '1.Características do Comportamento Empreendedor: Busca de oportunidades e iniciativa. Correr riscos calculados. Exigência de qualidade e eficiência. Persistência. Comprometimento. Busca de informações. Estabelecimento de metas. Monitoramento e planejamento sistemático. Persuasão e rede contatos. Independência e autoconfiança.2.Modelo de Negócios (Lean Canvas): Problema. Segmento de Clientes. Proposta de Valor Única. Solução. Métricas-Chave. Canais. Estrutura de Custos. Fluxos de Receita. Vantagem Injusta.3.Produto mínimo viável: Ciclo Construir-Mensurar-Aprender. Valor da vida útil do cliente.4.Plano de Negócios: Marketing, Finanças, Recursos Humanos, Desenvolvimento de Produtos e Tecnologia da Informação e Comunicação.'

# --- Row 17: Syllabus: label + EN full syllabus text (new B/C cells, style copied from row16) ---
$ws.Range("A17").Value = 'Syllabus:'
$ws.Range("B17").Value = '1.Characteristics of Entrepreneurial Behavior: Search for opportunities and initiative. Take calculated risks. Requirement of quality and efficiency. Persistence. Commitment. Search for information. Setting goals. Monitoring and systematic planning. Persuasion and network contacts. Independence and self-confidence.2.Business Model (Lean Canvas): Problem. Customer Segments. Unique Value Proposition. Solution. Key Metrics. Channels. Cost Structure. Revenue Streams. Unfair Advantage. 3.Minimum Viable Product: Build-Measure-Learn Cycle. Customer Lifetime Value.4.Business Plan: Marketing. Finance. Human Resources. Product Development. Information and communication technology.'
$ws.Range("C17").Value = '1.Characteristics of Entrepreneurial Behavior: Search for opportunities and initiative. Take calculated risks. Requirement of quality and efficiency. Persistence. Commitment. Search for information. Setting goals. Monitoring and systematic planning. Persuasion and network contacts. Independence and self-confidence.2.Business Model (Lean Canvas): Problem. Customer Segments. Unique Value Proposition. Solution. Key Metrics. Channels. Cost Structure. Revenue Streams. Unfair Advantage. 3.Minimum Viable Product: Build-Measure-Learn Cycle. Customer Lifetime Value.4.Business Plan: Marketing. Finance. Human Resources. Product Development. Information and communication technology.'
$ws.Range("B16").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("C16").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Rows.Item(17).RowHeight = 120

# --- Row 18: Avaliacao: label only; drop old B18/C18 ---
$ws.Range("A18").Value = 'Avaliação:'
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()
$ws.Rows.Item(18).AutoFit()

# --- Row 19: Metodo: label + method text (new position, value unchanged) ---
$ws.Range("A19").Value = 'Método:'

# --- Row 20: Criterio: label + evaluation text (new position, value unchanged) ---
$ws.Range("A20").Value = 'Critério:'

# --- Row 21: Norma de recuperacao: label + recovery text (new position, value unchanged); height 120 -> 60 ---
$ws.Range("A21").Value = 'Norma de recuperação:'
$ws.Rows.Item(21).RowHeight = 60

# --- Row 22: brand-new Bibliografia: row; copy styles from row 21 ---
$ws.Range("A22").Value = 'Bibliografia:'
$ws.Range("B22").Value = 'BLANK, Steve Gary. Do Sonho a realização em 4 passos: Estratégias para a criação de empresas de sucesso. Editora Evora. 3ª edição, 2008BLANK, Steve; DORF, Bob. STARTUP: Manual do Empreendedorismo. O guia passo a passo para construir uma grande empresa. Alta Books Editora.  1ª edição, 2014.CECCONELO, Antonio; AJZENTAL, Alberto. A construção do plano de negócios. Ed. Saraiva, 1ª edição, 2008.CHIAVENATO, Idalberto. Empreendedorismo – dando asas ao espírito empreendedor. Ed. Saraiva, 3ª edição, 2008.DOLABELA, Fernando. O Segredo de Luísa. Rio de Janeiro: Sextante, 2008. DORNELAS, Jose. Empreendedorismo: transformando ideias em negócios. Editora Campus. 1ª edição, 2001DORNELAS, Jose. Empreendedorismo na prática. LTC. 3ª edição, 2015DORNELAS, Jose Carlos Assis. Empreendedorismo na prática – mitos e verdades do empreendedor de sucesso. Elsevier/Campus: Rio de Janeiro, 2007. FILION, L. J.; Visão e Relações: Elementos para um Metamodelo da Atividade Empreendedora. International Small Business Journal, 1991. Tradução de Costa, S.R. FILION, L. J.; - O planejamento do seu Sistema de Aprendizagem Empresarial: Identifique uma Visão e Avalie o seu Sistema de Relações. Revista de Administração de Empresas, FGV, São Paulo, jul/set. 1991, pag. 31(3): 63:71. HASHIMOTO, Marcos. Espírito empreendedor nas organizações – aumentando a competitividade através do intraempreendedorismo. São Paulo: Saraiva, 2006. HISRICH, Robert; PETERS, Michael.  Empreendedorismo. 5.ed. - Porto Alegre: Bookman, 2004. OSTERWALDER, Alexander. Inovação Em Modelos de Negócios – Business Model Generation. Editora Alta Books, 2011PINCHOT, Gifford; PELLMAN, Ron. Intraempreendedorismo na prática: um guia de inovação. Campus: 2004RIES, Eric. A startup enxuta. Leya Editora. 1ª edição, 2011SANTOS. S.A. e CUNHA, N.C.V (orgs.). Empresas de Base Tecnológica: Conceitos, instrumentos e recursos. Unicorpore, 2005THIEL, Peter. De Zero a UM: O que aprender sobre empreendedorismo com Vale do Silício. Objetiva. 1ª edição, 2014TIMMONS; Jeffry; DORNELAS, José. SPINELLI, Stephen. A criação de novos negócios – empreendedorismo para o século 21. Editora Campus. 2010.'
$ws.Range("C22").Value = 'BLANK, Steve Gary. Do Sonho a realização em 4 passos: Estratégias para a criação de empresas de sucesso. Editora Evora. 3ª edição, 2008BLANK, Steve; DORF, Bob. STARTUP: Manual do Empreendedorismo. O guia passo a passo para construir uma grande empresa. Alta Books Editora.  1ª edição, 2014.CECCONELO, Antonio; AJZENTAL, Alberto. A construção do plano de negócios. Ed. Saraiva, 1ª edição, 2008.CHIAVENATO, Idalberto. Empreendedorismo – dando asas ao espírito empreendedor. Ed. Saraiva, 3ª edição, 2008.DOLABELA, Fernando. O Segredo de Luísa. Rio de Janeiro: Sextante, 2008. DORNELAS, Jose. Empreendedorismo: transformando ideias em negócios. Editora Campus. 1ª edição, 2001DORNELAS, Jose. Empreendedorismo na prática. LTC. 3ª edição, 2015DORNELAS, Jose Carlos Assis. Empreendedorismo na prática – mitos e verdades do empreendedor de sucesso. Elsevier/Campus: Rio de Janeiro, 2007. FILION, L. J.; Visão e Relações: Elementos para um Metamodelo da Atividade Empreendedora. International Small Business Journal, 1991. Tradução de Costa, S.R. FILION, L. J.; - O planejamento do seu Sistema de Aprendizagem Empresarial: Identifique uma Visão e Avalie o seu Sistema de Relações. Revista de Administração de Empresas, FGV, São Paulo, jul/set. 1991, pag. 31(3): 63:71. HASHIMOTO, Marcos. Espírito empreendedor nas organizações – aumentando a competitividade através do intraempreendedorismo. São Paulo: Saraiva, 2006. HISRICH, Robert; PETERS, Michael.  Empreendedorismo. 5.ed. - Porto Alegre: Bookman, 2004. OSTERWALDER, Alexander. Inovação Em Modelos de Negócios – Business Model Generation. Editora Alta Books, 2011PINCHOT, Gifford; PELLMAN, Ron. Intraempreendedorismo na prática: um guia de inovação. Campus: 2004RIES, Eric. A startup enxuta. Leya Editora. 1ª edição, 2011SANTOS. S.A. e CUNHA, N.C.V (orgs.). Empresas de Base Tecnológica: Conceitos, instrumentos e recursos. Unicorpore, 2005THIEL, Peter. De Zero a UM: O que aprender sobre empreendedorismo com Vale do Silício. Objetiva. 1ª edição, 2014TIMMONS; Jeffry; DORNELAS, José. SPINELLI, Stephen. A criação de novos negócios – empreendedorismo para o século 21. Editora Campus. 2010.'
$ws.Range("A21").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("B21").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("C21").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Rows.Item(22).RowHeight = 120

$excel.CutCopyMode = 0
